$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): insert two new columns ("height", "weight") between the
# existing "fumbles" (D1) and "fantasy points" (was E1) columns, pushing
# "fantasy points" out to the new G1 column.
# ---------------------------------------------------------------------------

# Clone the existing header style (bold/border/centered - style index 1) onto
# the two brand-new header cells by copying the already-styled E1 cell.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("E1").Copy($ws.Range("G1"))

# Move the "fantasy points" label to its new home (G1) before E1's text is
# overwritten, then fill in the new column headers.
$ws.Range("G1").Value = "fantasy points"
$ws.Range("F1").Value = "weight"
$ws.Range("E1").Value = "height"

# ---------------------------------------------------------------------------
# Data rows (2-17): every player-row gets a constant height/weight pair
# inserted in the new E/F columns, and its pre-existing "fantasy points"
# value (previously in column E) slides over to the new G column.
# ---------------------------------------------------------------------------
$height = 6.416666666666667
$weight = 255

for ($r = 2; $r -le 17; $r++) {
    $fantasyPoints = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 7).Value = $fantasyPoints
    $ws.Cells.Item($r, 5).Value = $height
    $ws.Cells.Item($r, 6).Value = $weight
}
